$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.326.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.24%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.308.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.86%  "

# Row 4
$ws.Range("E4").Value = "  -0.24%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.01%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.68%  "

# Row 7
$ws.Range("E7").Value = "  -0.22%  "

# Row 8
$ws.Range("E8").Value = "  -0.17%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.615"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.63%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.45"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.32%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0918"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.88%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.44"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.77%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.107"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.16%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.991"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.26%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.01%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.658.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.89%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.319.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.86%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.246.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.44%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.31%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000106"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.38%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.14%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.80%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "263.86"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.52%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.35"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.84%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.43%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.51%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.22%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.03%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "23.04"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.63%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.71"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.36%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "168.62"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.92%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0898"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.78%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.92"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.94%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.97"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.76%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.124"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +13.07%  "

# Row 36
$ws.Range("E36").Value = "  +1.08%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.68"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.87%  "

# Row 38
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.94"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +12.02%  "

# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0355"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.02%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.65"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.83%  "

# Row 41
$ws.Range("B41").Value = "BitcoinSV"
$ws.Range("C41").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "103.46"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +25.29%  "

# Row 42
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.51"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.39%  "

# Row 43
$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "71.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.10%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.230"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.09%  "

# Row 45
$ws.Range("E45").Value = "  +0.05%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.60"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +10.05%  "

# Row 47
$ws.Range("B47").Value = "ordi"
$ws.Range("C47").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "81.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +13.47%  "

# Row 48
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "114.81"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.09%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.25"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.85%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.37"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.50%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.28"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.70%  "
